$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text (matches the source data's
# inline-string type) while we overwrite their contents below; Excel would
# otherwise auto-coerce numeric-looking text into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '57.621.25'
$ws.Range('E2').Value = '  +1.95%  '

# Row 3
$ws.Range('D3').Value = '3.010.82'
$ws.Range('E3').Value = '  +0.58%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').Value = '510.30'
$ws.Range('E5').Value = '  +0.50%  '

# Row 6
$ws.Range('D6').Value = '139.36'
$ws.Range('E6').Value = '  +1.53%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').Value = '0.437'
$ws.Range('E8').Value = '  +1.34%  '

# Row 9
$ws.Range('D9').Value = '7.52'

# Row 10
$ws.Range('E10').Value = '  +1.68%  '

# Row 12
$ws.Range('D12').Value = '3.525.77'
$ws.Range('E12').Value = '  +0.51%  '

# Row 13
$ws.Range('E13').Value = '  +1.07%  '

# Row 14
$ws.Range('D14').Value = '26.45'

# Row 15
$ws.Range('E15').Value = '  +6.93%  '

# Row 16
$ws.Range('D16').Value = '57.588.75'
$ws.Range('E16').Value = '  +1.82%  '

# Row 17
$ws.Range('D17').Value = '6.24'
$ws.Range('E17').Value = '  +7.67%  '

# Row 18
$ws.Range('D18').Value = '3.014.49'
$ws.Range('E18').Value = '  +0.49%  '

# Row 19
$ws.Range('D19').Value = '12.82'
$ws.Range('E19').Value = '  +3.51%  '

# Row 20
$ws.Range('D20').Value = '7.97'
$ws.Range('E20').Value = '  +2.27%  '

# Row 21
$ws.Range('D21').Value = '331.08'
$ws.Range('E21').Value = '  +1.39%  '

# Row 22
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').Value = '0.498'
$ws.Range('E23').Value = '  +4.40%  '

# Row 24
$ws.Range('D24').Value = '64.53'
$ws.Range('E24').Value = '  +3.48%  '

# Row 25
$ws.Range('D25').Value = '0.170'
$ws.Range('E25').Value = '  -0.04%  '

# Row 26
$ws.Range('E26').Value = '  -0.09%  '

# Row 27
$ws.Range('E27').Value = '  +1.54%  '

# Row 28
$ws.Range('D28').Value = '6.77'
$ws.Range('E28').Value = '  +3.52%  '

# Row 29
$ws.Range('D29').Value = '7.38'
$ws.Range('E29').Value = '  +4.59%  '

# Row 30
$ws.Range('D30').Value = '1.82'
$ws.Range('E30').Value = '  +2.60%  '

# Row 31
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  -5.80%  '

# Row 32
$ws.Range('D32').Value = '20.60'
$ws.Range('E32').Value = '  +0.12%  '

# Row 33
$ws.Range('D33').Value = '4.71'
$ws.Range('E33').Value = '  +4.92%  '

# Row 34
$ws.Range('D34').Value = '153.92'
$ws.Range('E34').Value = '  -1.04%  '

# Row 35
$ws.Range('D35').Value = '5.85'
$ws.Range('E35').Value = '  +4.46%  '

# Row 36
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  +1.37%  '

# Row 37
$ws.Range('D37').Value = '24.40'
$ws.Range('E37').Value = '  +1.53%  '

# Row 38
$ws.Range('D38').Value = '0.0675'
$ws.Range('E38').Value = '  -0.02%  '

# Row 39
$ws.Range('D39').Value = '3.045.35'
$ws.Range('E39').Value = '  +0.54%  '

# Row 40
$ws.Range('D40').Value = '37.35'
$ws.Range('E40').Value = '  +1.29%  '

# Row 41
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.05%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').Value = '  +6.81%  '

# Row 43
$ws.Range('E43').Value = '  +0.60%  '

# Row 44
$ws.Range('D44').Value = '2.238.71'
$ws.Range('E44').Value = '  -0.92%  '

# Row 45
$ws.Range('E45').Value = '  +0.12%  '

# Row 46
$ws.Range('D46').Value = '0.984'
$ws.Range('E46').Value = '  -0.87%  '

# Row 47
$ws.Range('D47').Value = '6.01'
$ws.Range('E47').Value = '  +4.83%  '

# Row 48
$ws.Range('E48').Value = '  +1.74%  '

# Row 49
$ws.Range('D49').Value = '19.34'
$ws.Range('E49').Value = '  +1.64%  '

# Row 50
$ws.Range('E50').Value = '  -6.74%  '

# Row 51
$ws.Range('D51').Value = '0.0894'
$ws.Range('E51').Value = '  +2.69%  '

# Restore the original (default) cell style now that the text is committed,
# so formatting matches the source workbook (no lingering @ number format).
$ws.Range("D2:E51").Style = "Normal"